# Auto-generated script applying cell value updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 377.18182
$ws.Range("I11").Value = 377.18182
$ws.Range("K11").Value = 377.18182
$ws.Range("M11").Value = -237.18182
$ws.Range("H19").Value = 7645.5454
$ws.Range("J19").Value = 8901.111000000001
$ws.Range("L19").Value = 8901.111000000001
$ws.Range("N19").Value = -9251.111000000001
$ws.Range("H38").Value = 1526.6
$ws.Range("H41").Value = 133.75
$ws.Range("I41").Value = 225
$ws.Range("K41").Value = 225
$ws.Range("M41").Value = 215
$ws.Range("H62").Value = 4867.5557
$ws.Range("I62").Value = 5101
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 5101
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -4477
$ws.Range("N62").Value = -4248
$ws.Range("H63").Value = 68750
$ws.Range("J63").Value = 68750
$ws.Range("L63").Value = 68750
$ws.Range("N63").Value = -69998
$ws.Range("H64").Value = 6499.3335
$ws.Range("J64").Value = 6499.3335
$ws.Range("L64").Value = 6499.3335
$ws.Range("N64").Value = -6995.3335
$ws.Range("H65").Value = 4867.5557
$ws.Range("I65").Value = 5101
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 25505
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -22385
$ws.Range("N65").Value = -21240
$ws.Range("H66").Value = 68750
$ws.Range("J66").Value = 68750
$ws.Range("L66").Value = 206250
$ws.Range("N66").Value = -212490
$ws.Range("H67").Value = 6499.3335
$ws.Range("J67").Value = 6499.3335
$ws.Range("L67").Value = 6499.3335
$ws.Range("N67").Value = -8215.333500000001
$ws.Range("H70").Value = 3147.375
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 3454.1428
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 10362.4284
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -10902.4284
$ws.Range("H73").Value = 3147.375
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 3454.1428
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 10362.4284
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -12234.4284
$ws.Range("H74").Value = 5976.6924
$ws.Range("I74").Value = 5899
$ws.Range("K74").Value = 5899
$ws.Range("M74").Value = -4963
$ws.Range("H77").Value = 5976.6924
$ws.Range("I77").Value = 5899
$ws.Range("K77").Value = 29495
$ws.Range("M77").Value = -24815
$ws.Range("H82").Value = 374
$ws.Range("I82").Value = 374
$ws.Range("K82").Value = 1122
$ws.Range("M82").Value = -716
$ws.Range("H85").Value = 374
$ws.Range("I85").Value = 374
$ws.Range("K85").Value = 1122
$ws.Range("M85").Value = 282
$ws.Range("H98").Value = 15705.035
$ws.Range("I98").Value = 16266.346
$ws.Range("J98").Value = 269
$ws.Range("K98").Value = 16266.346
$ws.Range("L98").Value = 269
$ws.Range("M98").Value = -14768.346
$ws.Range("N98").Value = -3265
$ws.Range("H104").Value = 1323.75
$ws.Range("I104").Value = 1323.75
$ws.Range("K104").Value = 3971.25
$ws.Range("M104").Value = -2224.25
$ws.Range("H112").Value = 5265996
$ws.Range("J112").Value = 5684694.5
$ws.Range("L112").Value = 17054083.5
$ws.Range("N112").Value = -17056299.5
$ws.Range("H115").Value = 740
$ws.Range("I115").Value = 740
$ws.Range("K115").Value = 2220
$ws.Range("M115").Value = -653
$ws.Range("H122").Value = 15705.035
$ws.Range("I122").Value = 16266.346
$ws.Range("J122").Value = 269
$ws.Range("K122").Value = 48799.038
$ws.Range("L122").Value = 807
$ws.Range("M122").Value = -46349.038
$ws.Range("N122").Value = -5707
$ws.Range("H129").Value = 29020
$ws.Range("I129").Value = 7620.375
$ws.Range("K129").Value = 22861.125
$ws.Range("M129").Value = -17861.125
$ws.Range("H135").Value = 7814114
$ws.Range("I135").Value = 1100.591
$ws.Range("J135").Value = 25002744
$ws.Range("K135").Value = 9905.319
$ws.Range("L135").Value = 225024696
$ws.Range("M135").Value = -7370.319
$ws.Range("N135").Value = -225029766
$ws.Range("H137").Value = 48784950
$ws.Range("I137").Value = 38465108
$ws.Range("J137").Value = 66672680
$ws.Range("K137").Value = 115395324
$ws.Range("L137").Value = 200018040
$ws.Range("M137").Value = -115392774
$ws.Range("N137").Value = -200023140
$ws.Range("H138").Value = 3671834.8
$ws.Range("I138").Value = 4910.5
$ws.Range("J138").Value = 4025273.2
$ws.Range("K138").Value = 14731.5
$ws.Range("L138").Value = 12075819.6
$ws.Range("M138").Value = -9591.5
$ws.Range("N138").Value = -12086099.6
$ws.Range("H141").Value = 1695.5555
$ws.Range("I141").Value = 1629.7142
$ws.Range("K141").Value = 4889.142599999999
$ws.Range("M141").Value = 290.8574000000008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = -3226
$ws.Range("H5").Value = 271.75
$ws.Range("I5").Value = 231.66667
$ws.Range("K5").Value = 231.66667
$ws.Range("M5").Value = -119.66667
$ws.Range("H32").Value = 18527394
$ws.Range("I32").Value = 19238890
$ws.Range("J32").Value = 28457
$ws.Range("K32").Value = 19238890
$ws.Range("L32").Value = 28457
$ws.Range("M32").Value = -19238603
$ws.Range("N32").Value = -29031
$ws.Range("H45").Value = 2366.8
$ws.Range("I45").Value = 2369.8
$ws.Range("J45").Value = 2363.8
$ws.Range("K45").Value = 2369.8
$ws.Range("L45").Value = 2363.8
$ws.Range("M45").Value = -1992.8
$ws.Range("N45").Value = -3117.8
$ws.Range("H55").Value = 41666.668
$ws.Range("J55").Value = 55000
$ws.Range("L55").Value = 55000
$ws.Range("N55").Value = -55630
$ws.Range("H61").Value = 18522444
$ws.Range("I61").Value = 21742748
$ws.Range("J61").Value = 5686.625
$ws.Range("K61").Value = 21742748
$ws.Range("L61").Value = 5686.625
$ws.Range("M61").Value = -21742536
$ws.Range("N61").Value = -6110.625
$ws.Range("H74").Value = 51283788
$ws.Range("I74").Value = 86957900
$ws.Range("K74").Value = 86957900
$ws.Range("M74").Value = -86957026
$ws.Range("H77").Value = 51283788
$ws.Range("I77").Value = 86957900
$ws.Range("K77").Value = 434789500
$ws.Range("M77").Value = -434785132
$ws.Range("H102").Value = 63097.273
$ws.Range("I102").Value = 85053.71000000001
$ws.Range("J102").Value = 4546.778
$ws.Range("K102").Value = 85053.71000000001
$ws.Range("L102").Value = 4546.778
$ws.Range("M102").Value = -83431.71000000001
$ws.Range("N102").Value = -7790.778
$ws.Range("H110").Value = 53577
$ws.Range("I110").Value = 73008.8
$ws.Range("J110").Value = 4997.5
$ws.Range("K110").Value = 73008.8
$ws.Range("L110").Value = 4997.5
$ws.Range("M110").Value = -70963.8
$ws.Range("N110").Value = -9087.5
$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = $null
$ws.Range("N116").Value = -7588
$ws.Range("H122").Value = 3968.2307
$ws.Range("I122").Value = 3347.25
$ws.Range("J122").Value = 4244.222
$ws.Range("K122").Value = 10041.75
$ws.Range("L122").Value = 12732.666
$ws.Range("M122").Value = -7591.75
$ws.Range("N122").Value = -17632.666
$ws.Range("H132").Value = 31259530
$ws.Range("I132").Value = 10703.481
$ws.Range("J132").Value = 200003200
$ws.Range("K132").Value = 32110.443
$ws.Range("L132").Value = 600009600
$ws.Range("M132").Value = -29580.443
$ws.Range("N132").Value = -600014660
$ws.Range("H133").Value = 100906.6
$ws.Range("I133").Value = 89137
$ws.Range("J133").Value = 102214.336
$ws.Range("K133").Value = 89137
$ws.Range("L133").Value = 102214.336
$ws.Range("M133").Value = -86607
$ws.Range("N133").Value = -107274.336
$ws.Range("H135").Value = 45995.168
$ws.Range("J135").Value = 45995.168
$ws.Range("L135").Value = 45995.168
$ws.Range("N135").Value = -56135.168
$ws.Range("H136").Value = 18522444
$ws.Range("I136").Value = 21742748
$ws.Range("J136").Value = 5686.625
$ws.Range("K136").Value = 65228244
$ws.Range("L136").Value = 17059.875
$ws.Range("M136").Value = -65225694
$ws.Range("N136").Value = -22159.875
$ws.Range("H138").Value = 150000
$ws.Range("J138").Value = 150000
$ws.Range("L138").Value = 150000
$ws.Range("N138").Value = -160280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = -3228
$ws.Range("H4").Value = 271.75
$ws.Range("I4").Value = 231.66667
$ws.Range("K4").Value = 231.66667
$ws.Range("M4").Value = -116.66667
$ws.Range("H20").Value = 2512.9707
$ws.Range("I20").Value = 3153.9524
$ws.Range("J20").Value = 1477.5385
$ws.Range("K20").Value = 3153.9524
$ws.Range("L20").Value = 1477.5385
$ws.Range("M20").Value = -2906.9524
$ws.Range("N20").Value = -1971.5385
$ws.Range("H80").Value = 396.53845
$ws.Range("I80").Value = 577.2857
$ws.Range("J80").Value = 185.66667
$ws.Range("K80").Value = 577.2857
$ws.Range("L80").Value = 185.66667
$ws.Range("M80").Value = 420.7143
$ws.Range("N80").Value = -2181.66667
$ws.Range("H83").Value = 396.53845
$ws.Range("I83").Value = 577.2857
$ws.Range("J83").Value = 185.66667
$ws.Range("K83").Value = 2886.4285
$ws.Range("L83").Value = 928.3333500000001
$ws.Range("M83").Value = 2105.5715
$ws.Range("N83").Value = -10912.33335
$ws.Range("H86").Value = 36877.89
$ws.Range("I86").Value = 23181.2
$ws.Range("K86").Value = 23181.2
$ws.Range("M86").Value = -22058.2
$ws.Range("H89").Value = 36877.89
$ws.Range("I89").Value = 23181.2
$ws.Range("K89").Value = 115906
$ws.Range("M89").Value = -110290
$ws.Range("H94").Value = 2589.4375
$ws.Range("I94").Value = 2983.2
$ws.Range("J94").Value = 1933.1666
$ws.Range("K94").Value = 2983.2
$ws.Range("L94").Value = 1933.1666
$ws.Range("M94").Value = -2532.2
$ws.Range("N94").Value = -2835.1666
$ws.Range("H96").Value = 154994.5
$ws.Range("I96").Value = 19994
$ws.Range("K96").Value = 19994
$ws.Range("M96").Value = -17248
$ws.Range("H99").Value = 2734.8262
$ws.Range("I99").Value = 1461.25
$ws.Range("J99").Value = 5645.857
$ws.Range("K99").Value = 1461.25
$ws.Range("L99").Value = 5645.857
$ws.Range("M99").Value = 36.75
$ws.Range("N99").Value = -8641.857
$ws.Range("H105").Value = 14877.934
$ws.Range("I105").Value = 16889.924
$ws.Range("K105").Value = 16889.924
$ws.Range("M105").Value = -15142.924
$ws.Range("H107").Value = 6214.769
$ws.Range("I107").Value = 5663.364
$ws.Range("J107").Value = 9247.5
$ws.Range("K107").Value = 5663.364
$ws.Range("L107").Value = 9247.5
$ws.Range("M107").Value = -3743.364
$ws.Range("N107").Value = -13087.5
$ws.Range("H134").Value = 3290.9
$ws.Range("I134").Value = 3274.375
$ws.Range("J134").Value = 3357
$ws.Range("K134").Value = 9823.125
$ws.Range("L134").Value = 10071
$ws.Range("M134").Value = -7288.125
$ws.Range("N134").Value = -15141
$ws.Range("H141").Value = 88739.5
$ws.Range("J141").Value = 88739.5
$ws.Range("L141").Value = 88739.5
$ws.Range("N141").Value = -99099.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 702.9231
$ws.Range("I7").Value = 173.75
$ws.Range("K7").Value = 173.75
$ws.Range("M7").Value = -60.75
$ws.Range("H16").Value = 2185.1052
$ws.Range("I16").Value = 2317.375
$ws.Range("J16").Value = 1479.6666
$ws.Range("K16").Value = 2317.375
$ws.Range("L16").Value = 1479.6666
$ws.Range("M16").Value = -2030.375
$ws.Range("N16").Value = -2053.6666
$ws.Range("H18").Value = 46514.8
$ws.Range("J18").Value = 46514.8
$ws.Range("L18").Value = 46514.8
$ws.Range("N18").Value = -46974.8
$ws.Range("H22").Value = 7129.4
$ws.Range("I22").Value = 14500.857
$ws.Range("J22").Value = 679.375
$ws.Range("K22").Value = 14500.857
$ws.Range("L22").Value = 679.375
$ws.Range("M22").Value = -14150.857
$ws.Range("N22").Value = -1379.375
$ws.Range("H31").Value = 21744966
$ws.Range("I31").Value = 3898.923
$ws.Range("J31").Value = 50008350
$ws.Range("K31").Value = 3898.923
$ws.Range("L31").Value = 50008350
$ws.Range("M31").Value = -3603.923
$ws.Range("N31").Value = -50008940
$ws.Range("H34").Value = 21744966
$ws.Range("I34").Value = 3898.923
$ws.Range("J34").Value = 50008350
$ws.Range("K34").Value = 3898.923
$ws.Range("L34").Value = 50008350
$ws.Range("M34").Value = -3696.923
$ws.Range("N34").Value = -50008754
$ws.Range("H58").Value = 3335
$ws.Range("I58").Value = 3318.348
$ws.Range("J58").Value = 3462.6667
$ws.Range("K58").Value = 3318.348
$ws.Range("L58").Value = 3462.6667
$ws.Range("M58").Value = -3115.348
$ws.Range("N58").Value = -3868.6667
$ws.Range("H99").Value = 11345.083
$ws.Range("I99").Value = 15663.857
$ws.Range("J99").Value = 5298.8
$ws.Range("K99").Value = 15663.857
$ws.Range("L99").Value = 5298.8
$ws.Range("M99").Value = -14165.857
$ws.Range("N99").Value = -8294.799999999999
$ws.Range("H109").Value = 50001
$ws.Range("I109").Value = 50001
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 50001
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -48961
$ws.Range("N109").Value = $null
$ws.Range("H113").Value = 2185.1052
$ws.Range("I113").Value = 2317.375
$ws.Range("J113").Value = 1479.6666
$ws.Range("K113").Value = 2317.375
$ws.Range("L113").Value = 1479.6666
$ws.Range("M113").Value = -147.375
$ws.Range("N113").Value = -5819.6666
$ws.Range("H122").Value = 1792.3636
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 1871.6
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 5614.799999999999
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -10514.8
$ws.Range("H124").Value = 240024
$ws.Range("J124").Value = 240024
$ws.Range("L124").Value = 240024
$ws.Range("N124").Value = -244934
$ws.Range("H126").Value = 11345.083
$ws.Range("I126").Value = 15663.857
$ws.Range("J126").Value = 5298.8
$ws.Range("K126").Value = 46991.571
$ws.Range("L126").Value = 15896.4
$ws.Range("M126").Value = -44521.571
$ws.Range("N126").Value = -20836.4
$ws.Range("H132").Value = 147842.28
$ws.Range("I132").Value = 203775.1
$ws.Range("K132").Value = 611325.3
$ws.Range("M132").Value = -608795.3
$ws.Range("H134").Value = 1497.625
$ws.Range("I134").Value = 1395.15
$ws.Range("J134").Value = 2010
$ws.Range("K134").Value = 4185.450000000001
$ws.Range("L134").Value = 6030
$ws.Range("M134").Value = -1650.450000000001
$ws.Range("N134").Value = -11100
$ws.Range("H136").Value = 3335
$ws.Range("I136").Value = 3318.348
$ws.Range("J136").Value = 3462.6667
$ws.Range("K136").Value = 9955.044
$ws.Range("L136").Value = 10388.0001
$ws.Range("M136").Value = -7405.044
$ws.Range("N136").Value = -15488.0001
$ws.Range("H141").Value = 532963.0600000001
$ws.Range("I141").Value = 120000
$ws.Range("J141").Value = 564729.4399999999
$ws.Range("K141").Value = 120000
$ws.Range("L141").Value = 564729.4399999999
$ws.Range("M141").Value = -114820
$ws.Range("N141").Value = -575089.4399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 100000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 100000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 300000
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = -300540
$ws.Range("H67").Value = 100000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 100000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 300000
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = -301872
$ws.Range("H75").Value = 10083.333
$ws.Range("I75").Value = 12250
$ws.Range("J75").Value = 5750
$ws.Range("K75").Value = 36750
$ws.Range("L75").Value = 17250
$ws.Range("M75").Value = -35752
$ws.Range("N75").Value = -19246
$ws.Range("H78").Value = 10083.333
$ws.Range("I78").Value = 12250
$ws.Range("J78").Value = 5750
$ws.Range("K78").Value = 110250
$ws.Range("L78").Value = 51750
$ws.Range("M78").Value = -105258
$ws.Range("N78").Value = -61734
$ws.Range("H92").Value = 170
$ws.Range("J92").Value = 170
$ws.Range("L92").Value = 510
$ws.Range("N92").Value = -3006
$ws.Range("H113").Value = 5165.778
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5165.778
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 15497.334
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = -19837.334
$ws.Range("H122").Value = 3788.4167
$ws.Range("J122").Value = 4051.889
$ws.Range("L122").Value = 36467.001
$ws.Range("N122").Value = -41367.001
$ws.Range("H128").Value = 112938.6
$ws.Range("I128").Value = 112938.6
$ws.Range("K128").Value = 338815.8
$ws.Range("M128").Value = -333835.8
$ws.Range("H129").Value = 2511.6667
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 2511.6667
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 7535.000100000001
$ws.Range("M129").Value = $null
$ws.Range("N129").Value = -17535.0001
$ws.Range("H131").Value = 1845.3914
$ws.Range("I131").Value = 1845
$ws.Range("J131").Value = 1845.4736
$ws.Range("K131").Value = 5535
$ws.Range("L131").Value = 5536.4208
$ws.Range("M131").Value = -495
$ws.Range("N131").Value = -15616.4208
$ws.Range("H132").Value = 2567857.2
$ws.Range("J132").Value = 3513095
$ws.Range("L132").Value = 31617855
$ws.Range("N132").Value = -31622915
$ws.Range("H133").Value = 8013.9
$ws.Range("I133").Value = 6892.5
$ws.Range("J133").Value = 12499.5
$ws.Range("K133").Value = 20677.5
$ws.Range("L133").Value = 37498.5
$ws.Range("M133").Value = -15617.5
$ws.Range("N133").Value = -47618.5
$ws.Range("H134").Value = 7773
$ws.Range("I134").Value = 2254.4707
$ws.Range("J134").Value = 19499.875
$ws.Range("K134").Value = 6763.4121
$ws.Range("L134").Value = 58499.625
$ws.Range("M134").Value = -1693.4121
$ws.Range("N134").Value = -68639.625
$ws.Range("H137").Value = 3899.6667
$ws.Range("I137").Value = 3899.6667
$ws.Range("K137").Value = 11699.0001
$ws.Range("M137").Value = -6599.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3005380.8
$ws.Range("J3").Value = 6783.25
$ws.Range("L3").Value = 6783.25
$ws.Range("N3").Value = -7015.25
$ws.Range("H11").Value = 29547592
$ws.Range("I11").Value = 32501850
$ws.Range("K11").Value = 32501850
$ws.Range("M11").Value = -32501711
$ws.Range("H13").Value = 665.1111
$ws.Range("I13").Value = 167.33333
$ws.Range("K13").Value = 167.33333
$ws.Range("M13").Value = -28.33332999999999
$ws.Range("H70").Value = 122190.3
$ws.Range("I70").Value = 137872.67
$ws.Range("K70").Value = 137872.67
$ws.Range("M70").Value = -137602.67
$ws.Range("H73").Value = 122190.3
$ws.Range("I73").Value = 137872.67
$ws.Range("K73").Value = 137872.67
$ws.Range("M73").Value = -136936.67
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
$ws.Range("H113").Value = 4831.273
$ws.Range("I113").Value = 4763.4287
$ws.Range("J113").Value = 4950
$ws.Range("K113").Value = 4763.4287
$ws.Range("L113").Value = 4950
$ws.Range("M113").Value = -2593.4287
$ws.Range("N113").Value = -9290
$ws.Range("H122").Value = 3390.4443
$ws.Range("I122").Value = 2930.5715
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 8791.7145
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6341.7145
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 5344.405
$ws.Range("I132").Value = 4665.355
$ws.Range("J132").Value = 7258.091
$ws.Range("K132").Value = 13996.065
$ws.Range("L132").Value = 21774.273
$ws.Range("M132").Value = -11466.065
$ws.Range("N132").Value = -26834.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1174.8
$ws.Range("I16").Value = 1103.5385
$ws.Range("J16").Value = 1638
$ws.Range("K16").Value = 1103.5385
$ws.Range("L16").Value = 1638
$ws.Range("M16").Value = -933.5385000000001
$ws.Range("N16").Value = -1978
$ws.Range("H22").Value = 4331.353
$ws.Range("I22").Value = 3168.6
$ws.Range("K22").Value = 3168.6
$ws.Range("M22").Value = -2873.6
$ws.Range("H27").Value = 4331.353
$ws.Range("I27").Value = 3168.6
$ws.Range("K27").Value = 3168.6
$ws.Range("M27").Value = -3061.6
$ws.Range("H68").Value = 5129.5
$ws.Range("I68").Value = 3997.5
$ws.Range("J68").Value = 5506.8335
$ws.Range("K68").Value = 3997.5
$ws.Range("L68").Value = 5506.8335
$ws.Range("M68").Value = -3248.5
$ws.Range("N68").Value = -7004.8335
$ws.Range("H71").Value = 5129.5
$ws.Range("I71").Value = 3997.5
$ws.Range("J71").Value = 5506.8335
$ws.Range("K71").Value = 19987.5
$ws.Range("L71").Value = 27534.1675
$ws.Range("M71").Value = -16243.5
$ws.Range("N71").Value = -35022.1675
$ws.Range("H93").Value = 2963
$ws.Range("I93").Value = 1776.5
$ws.Range("J93").Value = 4149.5
$ws.Range("K93").Value = 1776.5
$ws.Range("L93").Value = 4149.5
$ws.Range("M93").Value = -528.5
$ws.Range("N93").Value = -6645.5
$ws.Range("H132").Value = 83338710
$ws.Range("I132").Value = 5112
$ws.Range("K132").Value = 15336
$ws.Range("M132").Value = -12806
$ws.Range("H136").Value = 7194.5
$ws.Range("I136").Value = 7194.5
$ws.Range("K136").Value = 21583.5
$ws.Range("M136").Value = -19033.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1749.5
$ws.Range("I6").Value = 1749.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1749.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1634.5
$ws.Range("N6").Value = $null
$ws.Range("H7").Value = 5000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = $null
$ws.Range("H100").Value = 76923550
$ws.Range("I100").Value = 100000470
$ws.Range("J100").Value = 478.66666
$ws.Range("K100").Value = 200000940
$ws.Range("L100").Value = 957.33332
$ws.Range("M100").Value = -200000399
$ws.Range("N100").Value = -2039.33332
$ws.Range("H107").Value = 1225
$ws.Range("I107").Value = 1225
$ws.Range("K107").Value = 3675
$ws.Range("M107").Value = -1755
$ws.Range("H113").Value = 791.46155
$ws.Range("I113").Value = 554.875
$ws.Range("J113").Value = 1170
$ws.Range("K113").Value = 1664.625
$ws.Range("L113").Value = 3510
$ws.Range("M113").Value = 505.375
$ws.Range("N113").Value = -7850
$ws.Range("H122").Value = 43525708
$ws.Range("I122").Value = 52687856
$ws.Range("K122").Value = 158063568
$ws.Range("M122").Value = -158061118
$ws.Range("H132").Value = 5970
$ws.Range("I132").Value = 6096.976
$ws.Range("K132").Value = 18290.928
$ws.Range("M132").Value = -15760.928
$ws.Range("H136").Value = 1217.3334
$ws.Range("I136").Value = 1216.5227
$ws.Range("J136").Value = 1226.25
$ws.Range("K136").Value = 3649.5681
$ws.Range("L136").Value = 3678.75
$ws.Range("M136").Value = -1099.5681
$ws.Range("N136").Value = -8778.75
